# Auto-generated edit script: update crypto price/volume table cells
# per the commit "Updated cryptos list on Fri Oct  6 23:47:57 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $newValue) {
    # Force the cell to remain a text value (matches the source data, which
    # stores numeric-looking prices like "213.42" / "1.00" as literal text),
    # then restore the default "Normal" style so no stray number format sticks.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '27.931.76'
Set-TextCell 'E2' '  +1.69%  '
Set-TextCell 'D3' '1.645.54'
Set-TextCell 'E3' '  +1.85%  '
Set-TextCell 'E4' '  -0.06%  '
Set-TextCell 'D5' '213.42'
Set-TextCell 'E5' '  +1.09%  '
Set-TextCell 'E6' '  +0.32%  '
Set-TextCell 'D7' '1.00'
Set-TextCell 'E7' '  -0.09%  '
Set-TextCell 'D8' '23.42'
Set-TextCell 'E8' '  +2.76%  '
Set-TextCell 'E9' '  +1.63%  '
Set-TextCell 'E10' '  +0.57%  '
Set-TextCell 'E11' '  -1.55%  '
Set-TextCell 'D12' '1.878.93'
Set-TextCell 'E12' '  +1.88%  '
Set-TextCell 'D13' '1.645.79'
Set-TextCell 'E13' '  +1.94%  '
Set-TextCell 'E14' '  +1.08%  '
Set-TextCell 'E15' '  +2.97%  '
Set-TextCell 'D16' '65.54'
Set-TextCell 'E16' '  +0.76%  '
Set-TextCell 'D17' '27.942.71'
Set-TextCell 'E17' '  +1.85%  '
Set-TextCell 'D18' '231.30'
Set-TextCell 'E18' '  -0.47%  '
Set-TextCell 'E19' '  +1.22%  '
Set-TextCell 'D20' '7.66'
Set-TextCell 'E20' '  +1.98%  '
Set-TextCell 'E21' '  -0.10%  '
Set-TextCell 'E22' '  +4.89%  '
Set-TextCell 'E23' '  +1.95%  '
Set-TextCell 'D24' '2.14'
Set-TextCell 'E24' '  +3.77%  '
Set-TextCell 'D25' '152.17'
Set-TextCell 'E25' '  +1.47%  '
Set-TextCell 'E27' '  +0.83%  '
Set-TextCell 'E28' '  +1.54%  '
Set-TextCell 'E29' '  -0.11%  '
Set-TextCell 'E30' '  +1.62%  '
Set-TextCell 'E31' '  +0.60%  '
Set-TextCell 'E32' '  +2.00%  '
Set-TextCell 'D33' '1.441.26'
Set-TextCell 'E33' '  -2.09%  '
Set-TextCell 'E34' '  +0.25%  '
Set-TextCell 'E35' '  +1.66%  '
Set-TextCell 'E36' '  -0.24%  '
Set-TextCell 'B37' 'TrustWalletToken'
Set-TextCell 'C37' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 'D37' '0.938'
Set-TextCell 'E37' '  -3.92%  '
Set-TextCell 'B38' 'ARBITRUM'
Set-TextCell 'C38' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 'D38' '0.890'
Set-TextCell 'E38' '  +3.33%  '
Set-TextCell 'E39' '  +1.25%  '
Set-TextCell 'D40' '0.560'
Set-TextCell 'E40' '  +0.98%  '
Set-TextCell 'D41' '69.09'
Set-TextCell 'E41' '  +3.37%  '
Set-TextCell 'E42' '  +3.65%  '
Set-TextCell 'E44' '  -0.26%  '
Set-TextCell 'D45' '1.83'
Set-TextCell 'E45' '  +6.00%  '
Set-TextCell 'E46' '  +3.68%  '
Set-TextCell 'E47' '  +0.72%  '
Set-TextCell 'D48' '1.787.61'
Set-TextCell 'E48' '  +1.60%  '
Set-TextCell 'D49' '89.05'
Set-TextCell 'E49' '  +2.75%  '
Set-TextCell 'E50' '  +0.23%  '
Set-TextCell 'E51' '  +1.14%  '
